$d = $word.ActiveDocument

# 1. "Requires Python 3 to run" -> "Requires Python 3.6+ to run"
$d.Content.Find.Execute("Requires Python 3 to run", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Requires Python 3.6+ to run", 2) | Out-Null

# 2. "Written by Chris Perkins in 2019:" -> "(c) 2019, Chris Perkins:"
$d.Content.Find.Execute("Written by Chris Perkins in 2019", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(c) 2019, Chris Perkins", 2) | Out-Null

# 3. Insert a new "v1.2 - code tidying." paragraph just before the
#    "v1.1 - fixed CSV output to UTF-8, fixed E.164 mask handling." paragraph.
$find = $d.Content
$found = $find.Find.Execute("v1.1 - fixed CSV output to UTF-8, fixed E.164 mask handling.")
if ($found) {
    $p = $find.Paragraphs(1)
    $r = $p.Range
    $r.InsertBefore("v1.2 - code tidying.`r")
}

# 4. Move the "_GoBack" bookmark from the Version History section up to the
#    very start of the document (start of the title paragraph), matching the
#    canonical location Word stores for the last edit position.
#    (Adding a bookmark with a name that already exists elsewhere relocates it,
#    removing the old bookmarkStart/bookmarkEnd pair automatically.)
$start = $d.Range(0, 0)
$start.InsertBefore("X")
$bmRange = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$delRange = $d.Range(0, 1)
$delRange.Delete()
